# SE-2885: Add example of GetHoldingsWithOrders
# Update the quotes.xlsx sample data: refresh the open/close price quotes
# for row 2 (AMZN) and move the active selection to H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update open_price / close_price values for the AMZN row
$ws.Range("F2").Value = 111
$ws.Range("G2").Value = 111.5

# Reflect the cursor/selection position saved with the workbook
$ws.Range("H2").Select()
